# Updated cryptos list with refreshed Price / Volume(1h) figures.
#
# Note: several "Price" values look like plain numbers (e.g. "195.83"), but
# Excel auto-converts such strings to numeric cells when assigned straight
# to .Value. To keep them as literal text (matching the sheet's original
# inlineStr/text cells), we temporarily switch the cell to Text number
# format ("@") before assigning the value, then restore the cell's original
# style/format afterwards so no visible formatting changes persist.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.536.66"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "3.501.94"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  +0.03%  "
$c_D5 = $ws.Range("D5")
$origStyle_D5 = $c_D5.Style
$c_D5.NumberFormat = "@"
$c_D5.Value = "602.75"
$c_D5.Style = $origStyle_D5
$ws.Range("E5").Value = "  -2.35%  "
$c_D6 = $ws.Range("D6")
$origStyle_D6 = $c_D6.Style
$c_D6.NumberFormat = "@"
$c_D6.Value = "195.83"
$c_D6.Style = $origStyle_D6
$ws.Range("E6").Value = "  +1.96%  "
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("E8").Value = "  +0.02%  "
$c_D9 = $ws.Range("D9")
$origStyle_D9 = $c_D9.Style
$c_D9.NumberFormat = "@"
$c_D9.Value = "0.199"
$c_D9.Style = $origStyle_D9
$ws.Range("E9").Value = "  -6.10%  "
$ws.Range("E10").Value = "  -1.85%  "
$c_D11 = $ws.Range("D11")
$origStyle_D11 = $c_D11.Style
$c_D11.NumberFormat = "@"
$c_D11.Value = "53.02"
$c_D11.Style = $origStyle_D11
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("E12").Value = "  -4.11%  "
$c_D13 = $ws.Range("D13")
$origStyle_D13 = $c_D13.Style
$c_D13.NumberFormat = "@"
$c_D13.Value = "9.41"
$c_D13.Style = $origStyle_D13
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("D14").Value = "4.060.82"
$ws.Range("E14").Value = "  -0.40%  "
$c_D15 = $ws.Range("D15")
$origStyle_D15 = $c_D15.Style
$c_D15.NumberFormat = "@"
$c_D15.Value = "597.83"
$c_D15.Style = $origStyle_D15
$ws.Range("E15").Value = "  -3.38%  "
$ws.Range("D16").Value = "69.717.93"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("E17").Value = "  -0.58%  "
$c_D18 = $ws.Range("D18")
$origStyle_D18 = $c_D18.Style
$c_D18.NumberFormat = "@"
$c_D18.Value = "12.60"
$c_D18.Style = $origStyle_D18
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("E19").Value = "  +2.50%  "
$ws.Range("D20").Value = "3.500.97"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("E21").Value = "  -0.90%  "
$c_D22 = $ws.Range("D22")
$origStyle_D22 = $c_D22.Style
$c_D22.NumberFormat = "@"
$c_D22.Value = "17.88"
$c_D22.Style = $origStyle_D22
$ws.Range("E22").Value = "  +4.25%  "
$ws.Range("E23").Value = "  +3.48%  "
$c_D24 = $ws.Range("D24")
$origStyle_D24 = $c_D24.Style
$c_D24.NumberFormat = "@"
$c_D24.Value = "101.56"
$c_D24.Style = $origStyle_D24
$ws.Range("E24").Value = "  -5.37%  "
$ws.Range("E25").Value = "  -2.34%  "
$c_D26 = $ws.Range("D26")
$origStyle_D26 = $c_D26.Style
$c_D26.NumberFormat = "@"
$c_D26.Value = "3.10"
$c_D26.Style = $origStyle_D26
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("E27").Value = "  -2.44%  "
$c_D28 = $ws.Range("D28")
$origStyle_D28 = $c_D28.Style
$c_D28.NumberFormat = "@"
$c_D28.Value = "9.43"
$c_D28.Style = $origStyle_D28
$ws.Range("E28").Value = "  -3.17%  "
$c_D29 = $ws.Range("D29")
$origStyle_D29 = $c_D29.Style
$c_D29.NumberFormat = "@"
$c_D29.Value = "32.90"
$c_D29.Style = $origStyle_D29
$ws.Range("E29").Value = "  -3.54%  "
$c_D30 = $ws.Range("D30")
$origStyle_D30 = $c_D30.Style
$c_D30.NumberFormat = "@"
$c_D30.Value = "4.27"
$c_D30.Style = $origStyle_D30
$ws.Range("E30").Value = "  +8.72%  "
$ws.Range("E31").Value = "  -0.44%  "
$c_D32 = $ws.Range("D32")
$origStyle_D32 = $c_D32.Style
$c_D32.NumberFormat = "@"
$c_D32.Value = "12.26"
$c_D32.Style = $origStyle_D32
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("E33").Value = "  -2.61%  "
$c_D34 = $ws.Range("D34")
$origStyle_D34 = $c_D34.Style
$c_D34.NumberFormat = "@"
$c_D34.Value = "63.16"
$c_D34.Style = $origStyle_D34
$ws.Range("E34").Value = "  -0.46%  "
$c_D35 = $ws.Range("D35")
$origStyle_D35 = $c_D35.Style
$c_D35.NumberFormat = "@"
$c_D35.Value = "3.16"
$c_D35.Style = $origStyle_D35
$ws.Range("E35").Value = "  +0.96%  "
$ws.Range("D36").Value = "3.719.12"
$ws.Range("E36").Value = "  +1.56%  "
$ws.Range("E38").Value = "  +2.07%  "
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("E40").Value = "  -1.99%  "
$c_D41 = $ws.Range("D41")
$origStyle_D41 = $c_D41.Style
$c_D41.NumberFormat = "@"
$c_D41.Value = "35.88"
$c_D41.Style = $origStyle_D41
$ws.Range("E41").Value = "  -2.33%  "
$c_D42 = $ws.Range("D42")
$origStyle_D42 = $c_D42.Style
$c_D42.NumberFormat = "@"
$c_D42.Value = "492.78"
$c_D42.Style = $origStyle_D42
$ws.Range("E42").Value = "  -4.31%  "
$c_D43 = $ws.Range("D43")
$origStyle_D43 = $c_D43.Style
$c_D43.NumberFormat = "@"
$c_D43.Value = "0.133"
$c_D43.Style = $origStyle_D43
$ws.Range("E43").Value = "  -3.66%  "
$c_D44 = $ws.Range("D44")
$origStyle_D44 = $c_D44.Style
$c_D44.NumberFormat = "@"
$c_D44.Value = "0.0449"
$c_D44.Style = $origStyle_D44
$ws.Range("E44").Value = "  -3.81%  "

# Rows 45 and 46: Stellar/ThetaToken swapped places, each also getting new Price/Volume values
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$c_D45 = $ws.Range("D45")
$origStyle_D45 = $c_D45.Style
$c_D45.NumberFormat = "@"
$c_D45.Value = "2.80"
$c_D45.Style = $origStyle_D45
$ws.Range("E45").Value = "  -4.30%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c_D46 = $ws.Range("D46")
$origStyle_D46 = $c_D46.Style
$c_D46.NumberFormat = "@"
$c_D46.Value = "0.138"
$c_D46.Style = $origStyle_D46
$ws.Range("E46").Value = "  -3.00%  "

$c_D47 = $ws.Range("D47")
$origStyle_D47 = $c_D47.Style
$c_D47.NumberFormat = "@"
$c_D47.Value = "3.24"
$c_D47.Style = $origStyle_D47
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("E49").Value = "  -4.36%  "
$ws.Range("E50").Value = "  +0.29%  "
$c_D51 = $ws.Range("D51")
$origStyle_D51 = $c_D51.Style
$c_D51.NumberFormat = "@"
$c_D51.Value = "128.25"
$c_D51.Style = $origStyle_D51
$ws.Range("E51").Value = "  -3.40%  "
